# Update the answer values in the "two-digit / one-digit division" table.
# Cells are addressed by (row, column) rather than by old text, because a
# couple of the new values coincide with old values elsewhere in the table
# (e.g. "29÷2=14, 1" is both an old value and a new value), which would make
# a simple global text Find/Replace ambiguous.
$d = $word.ActiveDocument
$t = $d.Tables(1)

$t.Cell(1, 1).Range.Text = "38÷2=19, 0"
$t.Cell(1, 2).Range.Text = "56÷9=6, 2"
$t.Cell(1, 3).Range.Text = "94÷6=15, 4"
$t.Cell(1, 4).Range.Text = "52÷6=8, 4"
$t.Cell(1, 5).Range.Text = "84÷5=16, 4"
$t.Cell(5, 1).Range.Text = "29÷2=14, 1"
$t.Cell(5, 2).Range.Text = "95÷2=47, 1"
$t.Cell(5, 3).Range.Text = "84÷2=42, 0"
$t.Cell(5, 4).Range.Text = "45÷6=7, 3"
$t.Cell(5, 5).Range.Text = "81÷5=16, 1"
$t.Cell(9, 1).Range.Text = "38÷6=6, 2"
$t.Cell(9, 2).Range.Text = "18÷8=2, 2"
$t.Cell(9, 3).Range.Text = "12÷8=1, 4"
$t.Cell(9, 4).Range.Text = "82÷2=41, 0"
$t.Cell(9, 5).Range.Text = "87÷4=21, 3"
$t.Cell(13, 1).Range.Text = "71÷9=7, 8"
$t.Cell(13, 2).Range.Text = "18÷5=3, 3"
$t.Cell(13, 3).Range.Text = "36÷9=4, 0"
$t.Cell(13, 4).Range.Text = "10÷3=3, 1"
$t.Cell(13, 5).Range.Text = "75÷3=25, 0"
$t.Cell(17, 1).Range.Text = "35÷4=8, 3"
$t.Cell(17, 2).Range.Text = "18÷9=2, 0"
$t.Cell(17, 3).Range.Text = "10÷7=1, 3"
$t.Cell(17, 4).Range.Text = "52÷4=13, 0"
$t.Cell(17, 5).Range.Text = "31÷3=10, 1"
